# "Generate Report for Handoff"
#
# The localization-status workbook tracks handoff state for each target
# language. This run marks the pending item as ready, stamps it with the
# generated handoff (.xlf) file + timestamp for zh-cn and de-de, and flips
# the dependency from "Ignored" to "Include" now that the handoff package
# exists.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handoff transform failed" -> "Ready for handoff" -------
# Update every occurrence so the shared string is fully replaced (not just
# appended) across all three sheets.
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsDeDe.Range("B2").Value = "Ready for handoff"

# --- zh-cn sheet: stamp the generated handoff file + datetime -------------
$zhFile = "67e9cb96-31c9-4dfe-acdd-11693bc2d527.0e87dcc5e31a544d0af291e79b2cd893763bb467.zh-cn.xlf"
$zhUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/4e752146205b18a6a7303d97b13d0331050dfd37/handoff/" + $zhFile

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $zhUrl, "", "", $zhFile) | Out-Null
$wsZhCn.Range("C2").Font.Underline = $true
$wsZhCn.Range("C2").Font.Color = 15570276

$wsZhCn.Range("D2").Value = "2016-02-26 07:00:59"
$wsZhCn.Range("H2").Value = "Include"

# --- de-de sheet: stamp the generated handoff file + datetime -------------
$deFile = "67e9cb96-31c9-4dfe-acdd-11693bc2d527.0e87dcc5e31a544d0af291e79b2cd893763bb467.de-de.xlf"
$deUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/4e752146205b18a6a7303d97b13d0331050dfd37/handoff/" + $deFile

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $deUrl, "", "", $deFile) | Out-Null
$wsDeDe.Range("C2").Font.Underline = $true
$wsDeDe.Range("C2").Font.Color = 15570276

$wsDeDe.Range("D2").Value = "2016-02-26 07:01:22"
$wsDeDe.Range("H2").Value = "Include"
